$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @(
    @{ Row = 2; Col = "D"; Value = 4.4506355143153 },
    @{ Row = 2; Col = "E"; Value = 0.00285604647398976 },
    @{ Row = 2; Col = "F"; Value = 9.02494084158365 },
    @{ Row = 3; Col = "D"; Value = -2.11609378373202 },
    @{ Row = 3; Col = "E"; Value = -12.8883611485623 },
    @{ Row = 3; Col = "F"; Value = 11.59381101243 },
    @{ Row = 4; Col = "D"; Value = -7.00930374180072 },
    @{ Row = 4; Col = "E"; Value = -16.6027367891976 },
    @{ Row = 4; Col = "F"; Value = 4.28883999589752 },
    @{ Row = 5; Col = "D"; Value = -1.74184002624819 },
    @{ Row = 5; Col = "E"; Value = -18.2670179550036 },
    @{ Row = 5; Col = "F"; Value = 22.7364232948649 },
    @{ Row = 6; Col = "D"; Value = -45.9509661291549 },
    @{ Row = 6; Col = "E"; Value = -52.2370300644955 },
    @{ Row = 6; Col = "F"; Value = -38.0131824436994 },
    @{ Row = 7; Col = "D"; Value = 20.2234112958771 },
    @{ Row = 7; Col = "E"; Value = 12.7463208602078 },
    @{ Row = 7; Col = "F"; Value = 28.108498415062 },
    @{ Row = 8; Col = "D"; Value = 21.2087247003941 },
    @{ Row = 8; Col = "E"; Value = 10.901569176131 },
    @{ Row = 8; Col = "F"; Value = 31.5510611421812 },
    @{ Row = 9; Col = "D"; Value = 57.9766805011587 },
    @{ Row = 9; Col = "E"; Value = 31.6400838022658 },
    @{ Row = 9; Col = "F"; Value = 95.9139171081479 },
    @{ Row = 10; Col = "D"; Value = 56.5810737412604 },
    @{ Row = 10; Col = "E"; Value = 30.3455798234838 },
    @{ Row = 10; Col = "F"; Value = 98.506311160116 },
    @{ Row = 11; Col = "D"; Value = -16.8657393725726 },
    @{ Row = 11; Col = "E"; Value = -25.1307165779118 },
    @{ Row = 11; Col = "F"; Value = -6.77511780061891 },
    @{ Row = 12; Col = "D"; Value = 3.43294638693178 },
    @{ Row = 12; Col = "E"; Value = -0.65012984413163 },
    @{ Row = 12; Col = "F"; Value = 7.76134619583574 },
    @{ Row = 13; Col = "D"; Value = 19.8472512305629 },
    @{ Row = 13; Col = "E"; Value = 6.69157486838403 },
    @{ Row = 13; Col = "F"; Value = 34.1299304695193 },
    @{ Row = 14; Col = "D"; Value = 120.130500317897 },
    @{ Row = 14; Col = "E"; Value = 63.7314956449412 },
    @{ Row = 14; Col = "F"; Value = 214.822381699593 },
    @{ Row = 15; Col = "D"; Value = 210.557828156942 },
    @{ Row = 15; Col = "E"; Value = 97.0460216264289 },
    @{ Row = 15; Col = "F"; Value = 398.826930346321 },
    @{ Row = 16; Col = "D"; Value = 145.771235644729 },
    @{ Row = 16; Col = "E"; Value = 83.5447401488224 },
    @{ Row = 16; Col = "F"; Value = 230.641682750723 },
    @{ Row = 17; Col = "D"; Value = 4.07448859494469 },
    @{ Row = 17; Col = "E"; Value = 0.839500229598382 },
    @{ Row = 17; Col = "F"; Value = 7.50010627375154 },
    @{ Row = 18; Col = "D"; Value = 21.9194187208039 },
    @{ Row = 18; Col = "E"; Value = 13.9518720111582 },
    @{ Row = 18; Col = "F"; Value = 31.3962266751789 },
    @{ Row = 19; Col = "D"; Value = 76.7214989760613 },
    @{ Row = 19; Col = "E"; Value = 53.4462798150482 },
    @{ Row = 19; Col = "F"; Value = 105.219584370388 },
    @{ Row = 20; Col = "D"; Value = 66.9266910978088 },
    @{ Row = 20; Col = "E"; Value = 31.2964152129587 },
    @{ Row = 20; Col = "F"; Value = 117.638330980901 },
    @{ Row = 21; Col = "D"; Value = -43.4121490548788 },
    @{ Row = 21; Col = "E"; Value = -48.6899213714106 },
    @{ Row = 21; Col = "F"; Value = -36.7459298937381 }
)

foreach ($u in $updates) {
    $ws.Range("$($u.Col)$($u.Row)").Value = $u.Value
}
